$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row before row 11. This shifts the old rows 11-14 down to 12-15
#    and automatically grows the merged cell B6:B12 -> B6:B13 and moves B13:B14 -> B14:B15.
$ws.Rows.Item(11).Insert()

# 2. The newly inserted row 11 has no formatting yet; copy the formatting from the
#    row above (row 8), which uses the same visual pattern (gray B cell, bordered C/D cells).
$ws.Range("B8").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the row height used by the other data rows in this table.
$ws.Rows.Item(11).RowHeight = 19.5

# 3. New content in the freshly inserted row: "5. Verifica adição de stock"
$ws.Range("D11").Value = "5. Verifica adição de stock"

# 4. Renumber the two steps that used to be 5 and 6, which are now rows 12 and 13.
$ws.Range("D12").Value = "6. Regista adição do stock"
$ws.Range("D13").Value = "7. Indica que a adição  foi feita com sucesso"

# 5. Row 14 (previously row 13, the exception row) gets new wording, and the old
#    exception text moves down to the new trailing row 15.
$ws.Range("D15").Value = "4.1. Apresenta mensagem ""Stock Inválido"""
$ws.Range("D14").Value = "4.1 Dados da peça inválidos"

# 6. Keep the active selection in sync with the new layout.
$ws.Range("E11").Select()
